# "emploi du temps" workbook update
# - row 14..18 (B column) task description gets extended with the new
#   research items (shared string text change)
# - column B is widened to fit the new (longer) text
# - selection / scroll position moves to the area being edited (B14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "analyse, recherche et test de différents algorithmes ML + recherche liste de mots en anglais par émotions + recherche d'api"

# Update the 5 task cells (31 Jan - 4 Feb) with the extended description.
$ws.Range("B14:B18").Value = $newText

# Widen column B so the longer text keeps fitting (mirrors the author's
# manual "best fit" width bump from 63 to ~105.7 characters).
$ws.Columns.Item(2).ColumnWidth = 104.8776041666667

# Move the selection/viewport to where the edit happened.
$ws.Range("B14").Select()
